$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.1948051948051948
$ws.Range("C2").Value = 0.5324675324675324
$ws.Range("J2").Value = 0.02337662337662338
$ws.Range("P2").Value = 0.1298701298701299
$ws.Range("S2").Value = 0.1194805194805195

# Row 3
$ws.Range("B3").Value = 0.0091324200913242
$ws.Range("C3").Value = 0.0547945205479452
$ws.Range("J3").Value = 0.0547945205479452
$ws.Range("P3").Value = 0.6940639269406392
$ws.Range("S3").Value = 0.1872146118721461

# Row 4
$ws.Range("J4").Value = 0.0576923076923077
$ws.Range("P4").Value = 0.7115384615384616
$ws.Range("S4").Value = 0.2307692307692308

# Row 5
$ws.Range("P5").Value = 0.4
$ws.Range("S5").Value = 0.6

# Row 6
$ws.Range("B6").Value = 0.06
$ws.Range("D6").Value = 0.02
$ws.Range("F6").Value = 0.044
$ws.Range("J6").Value = 0.248
$ws.Range("O6").Value = 0.02
$ws.Range("Q6").Value = 0.2
$ws.Range("R6").Value = 0.06
$ws.Range("S6").Value = 0.348

# Row 7
$ws.Range("B7").Value = 0.1359649122807018
$ws.Range("D7").Value = 0.02192982456140351
$ws.Range("E7").Value = 0.004385964912280702
$ws.Range("F7").Value = 0.07017543859649122
$ws.Range("J7").Value = 0.1271929824561404
$ws.Range("O7").Value = 0.0131578947368421
$ws.Range("Q7").Value = 0.1535087719298246
$ws.Range("R7").Value = 0.08771929824561403
$ws.Range("S7").Value = 0.3859649122807017

# Row 8
$ws.Range("B8").Value = 0.1085814360770578
$ws.Range("D8").Value = 0.01751313485113835
$ws.Range("E8").Value = 0.001751313485113835
$ws.Range("F8").Value = 0.08056042031523643
$ws.Range("J8").Value = 0.1348511383537653
$ws.Range("O8").Value = 0.01751313485113835
$ws.Range("Q8").Value = 0.2031523642732049
$ws.Range("R8").Value = 0.05779334500875657
$ws.Range("S8").Value = 0.3782837127845884

# Row 9
$ws.Range("B9").Value = 0.1161616161616162
$ws.Range("D9").Value = 0.0505050505050505
$ws.Range("F9").Value = 0.09090909090909091
$ws.Range("J9").Value = 0.1616161616161616
$ws.Range("O9").Value = 0.01515151515151515
$ws.Range("Q9").Value = 0.1767676767676768
$ws.Range("R9").Value = 0.0505050505050505
$ws.Range("S9").Value = 0.3383838383838384

# Row 10
$ws.Range("B10").Value = 0.1186079545454545
$ws.Range("D10").Value = 0.01704545454545454
$ws.Range("E10").Value = 0.002130681818181818
$ws.Range("F10").Value = 0.06534090909090909
$ws.Range("J10").Value = 0.1292613636363636
$ws.Range("O10").Value = 0.01633522727272727
$ws.Range("Q10").Value = 0.2052556818181818
$ws.Range("R10").Value = 0.08522727272727272
$ws.Range("S10").Value = 0.3607954545454545

# Row 11
$ws.Range("F11").Value = 0.002906976744186046
$ws.Range("G11").Value = 0.1511627906976744
$ws.Range("J11").Value = 0.0872093023255814
$ws.Range("K11").Value = 0.2383720930232558
$ws.Range("L11").Value = 0.4941860465116279
$ws.Range("S11").Value = 0.02616279069767442

# Row 12
$ws.Range("G12").Value = 0.7758620689655172
$ws.Range("J12").Value = 0.1609195402298851
$ws.Range("K12").Value = 0.005747126436781609
$ws.Range("L12").Value = 0.01149425287356322
$ws.Range("S12").Value = 0.04597701149425287

# Row 13
$ws.Range("G13").Value = 0.65
$ws.Range("J13").Value = 0.275
$ws.Range("S13").Value = 0.075

# Row 15
$ws.Range("F15").Value = 0.009302325581395349
$ws.Range("H15").Value = 0.1581395348837209
$ws.Range("I15").Value = 0.09302325581395349
$ws.Range("J15").Value = 0.3302325581395349
$ws.Range("K15").Value = 0.08372093023255814
$ws.Range("O15").Value = 0.05581395348837209
$ws.Range("S15").Value = 0.2697674418604651

# Row 16
$ws.Range("F16").Value = 0.02164502164502164
$ws.Range("H16").Value = 0.1601731601731602
$ws.Range("I16").Value = 0.06926406926406926
$ws.Range("J16").Value = 0.3246753246753247
$ws.Range("K16").Value = 0.1212121212121212
$ws.Range("M16").Value = 0.0735930735930736
$ws.Range("O16").Value = 0.0735930735930736
$ws.Range("S16").Value = 0.1558441558441558

# Row 17
$ws.Range("F17").Value = 0.01923076923076923
$ws.Range("H17").Value = 0.1942307692307692
$ws.Range("I17").Value = 0.09230769230769231
$ws.Range("J17").Value = 0.425
$ws.Range("K17").Value = 0.08076923076923077
$ws.Range("M17").Value = 0.02115384615384616
$ws.Range("N17").Value = 0.001923076923076923
$ws.Range("O17").Value = 0.03076923076923077
$ws.Range("S17").Value = 0.1346153846153846

# Row 18
$ws.Range("F18").Value = 0.0303030303030303
$ws.Range("H18").Value = 0.2222222222222222
$ws.Range("I18").Value = 0.08080808080808081
$ws.Range("J18").Value = 0.3636363636363636
$ws.Range("K18").Value = 0.07575757575757576
$ws.Range("M18").Value = 0.03535353535353535
$ws.Range("N18").Value = 0.005050505050505051
$ws.Range("O18").Value = 0.05555555555555555
$ws.Range("S18").Value = 0.1313131313131313

# Row 19
$ws.Range("F19").Value = 0.009478672985781991
$ws.Range("H19").Value = 0.2417061611374408
$ws.Range("I19").Value = 0.06635071090047394
$ws.Range("J19").Value = 0.3405551794177387
$ws.Range("K19").Value = 0.1056194989844279
$ws.Range("M19").Value = 0.02979011509817197
$ws.Range("N19").Value = 0.0006770480704129993
$ws.Range("O19").Value = 0.05619498984427895
$ws.Range("S19").Value = 0.1496276235612728
